$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset's merge/link columns were renamed as part of the
# "create core Dataset object" refactor:
#   PIDN_link     -> PIDN_x
#   DCDate_link   -> DCDate_x
#   InstrID_link  -> InstrID_x
#   _merge        -> _mp_merge
#   _diff_days    -> _mp_diff_days
#   _duplicates   -> _mp_duplicates
$replacements = @{
    "PIDN_link"    = "PIDN_x"
    "DCDate_link"  = "DCDate_x"
    "InstrID_link" = "InstrID_x"
    "_merge"       = "_mp_merge"
    "_diff_days"   = "_mp_diff_days"
    "_duplicates"  = "_mp_duplicates"
}

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($replacements.ContainsKey($val)) {
        $cell.Value = $replacements[$val]
    }
}

# Column headers changed length, so re-fit the column widths to the
# new (longer) header text, same as Excel does automatically when the
# sheet is re-saved after editing cell content.
$used.EntireColumn.AutoFit()
